$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 (sub_id 106) ---
$ws.Cells.Item(6, 1).Value = 45972
$ws.Cells.Item(6, 2).Value = 106
$ws.Cells.Item(6, 4).Value = "R_1kNcuYDMTmIPVBg"

# --- Row 8 (sub_id 108) ---
$ws.Cells.Item(8, 1).Value = 45972
$ws.Cells.Item(8, 2).Value = 108
$ws.Cells.Item(8, 4).Value = "R_325Icq7AOez25LX"

# --- back to row 6 column C ---
$ws.Cells.Item(6, 3).Value = "R_3VL7UmQBWN0n2Xn"
$ws.Cells.Item(6, 5).Value = "y"

# --- row 7 (sub_id 107) ---
$ws.Cells.Item(7, 2).Value = 107
$ws.Cells.Item(7, 5).Value = "y"

# --- row 8 column C ---
$ws.Cells.Item(8, 3).Value = "R_5NlSb3U49Wgqq9j"
$ws.Cells.Item(8, 5).Value = "n"

# --- row 9 (sub_id 109) ---
$ws.Cells.Item(9, 2).Value = 109
$ws.Cells.Item(9, 5).Value = "y"

# Match the date formatting used by the existing rows (style used for column A)
$ws.Range("A2").Copy()
$ws.Range("A6:A9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("E12").Select()
